$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) columns
# D-column values are forced to Text format so numeric-looking
# strings (e.g. "20.51") are stored verbatim instead of being
# auto-converted to floating point numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.255.56"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.660.25"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.69%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.47"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5327"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2637"
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06347"
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.51"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.541"
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.692.64"
$ws.Range("E13").Value = "  +1.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.888.16"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5518"
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8179"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.252.49"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.649"
$ws.Range("E20").Value = "  +2.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.02"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.052"
$ws.Range("E23").Value = "  +1.28%  "
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.65"
$ws.Range("E25").Value = "  +3.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1228"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.06"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.473"
$ws.Range("E29").Value = "  +2.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05799"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.278"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.572"
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.286"
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.612"
$ws.Range("E34").Value = "  +4.02%  "
$ws.Range("E35").Value = "  +2.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9581"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.430"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5795"
$ws.Range("E38").Value = "  +2.65%  "
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.850"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8516"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.74"
$ws.Range("E43").Value = "  +3.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.046.86"
$ws.Range("E44").Value = "  +3.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.801.37"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.10"
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₈106"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.013"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4369"
$ws.Range("E49").Value = "  +1.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.967"
$ws.Range("E50").Value = "  +2.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05158"
$ws.Range("E51").Value = "  +0.17%  "
